# Updated cryptos list with refreshed prices / 1h volume changes.
# D-column values are prefixed with a leading apostrophe so Excel stores
# them as text (matching the original inlineStr cells) instead of
# re-interpreting number-like strings (e.g. "1.000", "27.161.10").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.161.10'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').Value = '''1.871.70'
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''307.32'
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '''0.5119'
$ws.Range('E7').Value = '  +2.29%  '
$ws.Range('D8').Value = '''0.3747'
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('D9').Value = '''0.07138'
$ws.Range('E9').Value = '  -2.08%  '
$ws.Range('D10').Value = '''0.8871'
$ws.Range('E10').Value = '  -2.61%  '
$ws.Range('D11').Value = '''20.62'
$ws.Range('E11').Value = '  -3.27%  '
$ws.Range('D12').Value = '''1.864.98'
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('D14').Value = '''5.325'
$ws.Range('E14').Value = '  -2.71%  '
$ws.Range('D15').Value = '''89.18'
$ws.Range('E15').Value = '  -3.62%  '
$ws.Range('D16').Value = '''1.001'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').Value = '''0.000008469'
$ws.Range('E17').Value = '  -2.99%  '
$ws.Range('D18').Value = '''14.08'
$ws.Range('E18').Value = '  -4.03%  '
$ws.Range('D19').Value = '''1.000'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = '''27.200.33'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').Value = '''5.055'
$ws.Range('E21').Value = '  -2.25%  '
$ws.Range('D22').Value = '''2.107.70'
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('E23').Value = '  -2.98%  '
$ws.Range('D24').Value = '''6.477'
$ws.Range('E24').Value = '  -1.89%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''149.96'
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '''1.848'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '''17.96'
$ws.Range('D28').Value = '''2.098'
$ws.Range('E28').Value = '  -5.61%  '
$ws.Range('D29').Value = '''112.89'
$ws.Range('E29').Value = '  -1.94%  '
$ws.Range('D30').Value = '''4.716'
$ws.Range('E30').Value = '  -3.60%  '
$ws.Range('D31').Value = '''4.668'
$ws.Range('E31').Value = '  -3.35%  '
$ws.Range('D32').Value = '''0.09024'
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('D33').Value = '''0.05136'
$ws.Range('E33').Value = '  -2.87%  '
$ws.Range('D34').Value = '''3.086'
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('D35').Value = '''1.157'
$ws.Range('E35').Value = '  -6.52%  '
$ws.Range('D36').Value = '''0.7349'
$ws.Range('E36').Value = '  -6.71%  '
$ws.Range('D37').Value = '''0.02047'
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('D38').Value = '''2.508'
$ws.Range('E38').Value = '  -5.91%  '
$ws.Range('D39').Value = '''3.058'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').Value = '''1.075'
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('D41').Value = '''0.5364'
$ws.Range('E41').Value = '  -3.13%  '
$ws.Range('D42').Value = '''6.581'
$ws.Range('E42').Value = '  -3.25%  '
$ws.Range('D43').Value = '''117.12'
$ws.Range('E43').Value = '  +2.83%  '
$ws.Range('D44').Value = '''8.324'
$ws.Range('E44').Value = '  -2.48%  '
$ws.Range('D45').Value = '''0.1472'
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range('D46').Value = '''0.4640'
$ws.Range('E46').Value = '  -3.83%  '
$ws.Range('D47').Value = '''0.9999'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').Value = '''10.05'
$ws.Range('E48').Value = '  -5.35%  '
$ws.Range('E49').Value = '  -4.11%  '
$ws.Range('D50').Value = '''64.39'
$ws.Range('E50').Value = '  -4.38%  '
$ws.Range('E51').Value = '  -1.38%  '
